# Update the "Spices" food_group rows: set type (column F) from "B" to "A"
# and categoryY2013 (column G) to 5, reflecting the second reference diet
# results that were added and recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1731
$endRow = 1947

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"
    $ws.Cells.Item($r, 7).Value = 5
}
